# Apply updated timing results (redid data with new sizes)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Brute Force")
$ws.Range("B2").Value = 0.000003576278686523438
$ws.Range("C2").Value = 0.000003814697265625
$ws.Range("D2").Value = 0.000001192092895507812
$ws.Range("E2").Value = 0.000002145767211914062
$ws.Range("F2").Value = 0.0000019073486328125
$ws.Range("G2").Value = 0.000002145767211914062
$ws.Range("H2").Value = 0.000001192092895507812
$ws.Range("I2").Value = 0.0000007152557373046875
$ws.Range("J2").Value = 0.0000007152557373046875
$ws.Range("K2").Value = 0.000001192092895507812
$ws.Range("A3").Value = 200
$ws.Range("B3").Value = 0.04208827018737793
$ws.Range("C3").Value = 0.04571247100830078
$ws.Range("D3").Value = 0.0386815071105957
$ws.Range("E3").Value = 0.04149127006530762
$ws.Range("F3").Value = 0.0389859676361084
$ws.Range("G3").Value = 0.04320383071899414
$ws.Range("H3").Value = 0.03923869132995605
$ws.Range("I3").Value = 0.03976917266845703
$ws.Range("J3").Value = 0.04477643966674805
$ws.Range("K3").Value = 0.04109621047973633
$ws.Range("A4").Value = 400
$ws.Range("B4").Value = 0.1566979885101318
$ws.Range("C4").Value = 0.160341739654541
$ws.Range("D4").Value = 0.1572701930999756
$ws.Range("E4").Value = 0.1587464809417725
$ws.Range("F4").Value = 0.1577141284942627
$ws.Range("G4").Value = 0.1605966091156006
$ws.Range("H4").Value = 0.1619212627410889
$ws.Range("I4").Value = 0.1600699424743652
$ws.Range("J4").Value = 0.1562099456787109
$ws.Range("K4").Value = 0.160663366317749
$ws.Range("A5").Value = 600
$ws.Range("B5").Value = 0.3513696193695068
$ws.Range("C5").Value = 0.3570907115936279
$ws.Range("D5").Value = 0.3501062393188477
$ws.Range("E5").Value = 0.3559532165527344
$ws.Range("F5").Value = 0.3539128303527832
$ws.Range("G5").Value = 0.3547420501708984
$ws.Range("H5").Value = 0.3518776893615723
$ws.Range("I5").Value = 0.3583123683929443
$ws.Range("J5").Value = 0.3493270874023438
$ws.Range("K5").Value = 0.3555731773376465
$ws.Range("A6").Value = 800
$ws.Range("B6").Value = 0.6416592597961426
$ws.Range("C6").Value = 0.6862530708312988
$ws.Range("D6").Value = 0.6367764472961426
$ws.Range("E6").Value = 0.7484517097473145
$ws.Range("F6").Value = 0.6362960338592529
$ws.Range("G6").Value = 0.6355326175689697
$ws.Range("H6").Value = 0.6237063407897949
$ws.Range("I6").Value = 0.6676595211029053
$ws.Range("J6").Value = 0.6421687602996826
$ws.Range("K6").Value = 0.6449816226959229
$ws.Range("A7").Value = 1000
$ws.Range("B7").Value = 1.028444051742554
$ws.Range("C7").Value = 1.009976863861084
$ws.Range("D7").Value = 1.027829647064209
$ws.Range("E7").Value = 1.065417528152466
$ws.Range("F7").Value = 1.044753789901733
$ws.Range("G7").Value = 1.029937744140625
$ws.Range("H7").Value = 0.9980900287628174
$ws.Range("I7").Value = 1.02228856086731
$ws.Range("J7").Value = 1.019595861434937
$ws.Range("K7").Value = 1.040081262588501

$ws = $wb.Worksheets.Item("Divide and Conquer")
$ws.Range("B2").Value = 0.00002980232238769531
$ws.Range("C2").Value = 0.0000133514404296875
$ws.Range("D2").Value = 0.00000476837158203125
$ws.Range("E2").Value = 0.000008344650268554688
$ws.Range("F2").Value = 0.000006198883056640625
$ws.Range("G2").Value = 0.000008821487426757812
$ws.Range("H2").Value = 0.000004291534423828125
$ws.Range("I2").Value = 0.000003099441528320312
$ws.Range("J2").Value = 0.00000286102294921875
$ws.Range("K2").Value = 0.00000286102294921875
$ws.Range("A3").Value = 200
$ws.Range("B3").Value = 0.001081705093383789
$ws.Range("C3").Value = 0.001033782958984375
$ws.Range("D3").Value = 0.0009775161743164062
$ws.Range("E3").Value = 0.001008749008178711
$ws.Range("F3").Value = 0.001081228256225586
$ws.Range("G3").Value = 0.001187562942504883
$ws.Range("H3").Value = 0.001106739044189453
$ws.Range("I3").Value = 0.001053094863891602
$ws.Range("J3").Value = 0.001139402389526367
$ws.Range("K3").Value = 0.001062154769897461
$ws.Range("A4").Value = 400
$ws.Range("B4").Value = 0.002152442932128906
$ws.Range("C4").Value = 0.002050399780273438
$ws.Range("D4").Value = 0.00194096565246582
$ws.Range("E4").Value = 0.001989126205444336
$ws.Range("F4").Value = 0.001961231231689453
$ws.Range("G4").Value = 0.002144575119018555
$ws.Range("H4").Value = 0.002060413360595703
$ws.Range("I4").Value = 0.002064228057861328
$ws.Range("J4").Value = 0.002092599868774414
$ws.Range("K4").Value = 0.001862525939941406
$ws.Range("A5").Value = 600
$ws.Range("B5").Value = 0.002901315689086914
$ws.Range("C5").Value = 0.003487348556518555
$ws.Range("D5").Value = 0.003481626510620117
$ws.Range("E5").Value = 0.003158330917358398
$ws.Range("F5").Value = 0.002799510955810547
$ws.Range("G5").Value = 0.002933502197265625
$ws.Range("H5").Value = 0.002854108810424805
$ws.Range("I5").Value = 0.002927303314208984
$ws.Range("J5").Value = 0.002887725830078125
$ws.Range("K5").Value = 0.002881526947021484
$ws.Range("A6").Value = 800
$ws.Range("B6").Value = 0.004082441329956055
$ws.Range("C6").Value = 0.003956794738769531
$ws.Range("D6").Value = 0.008474349975585938
$ws.Range("E6").Value = 0.004492044448852539
$ws.Range("F6").Value = 0.004628181457519531
$ws.Range("G6").Value = 0.003955841064453125
$ws.Range("H6").Value = 0.003919601440429688
$ws.Range("I6").Value = 0.004063129425048828
$ws.Range("J6").Value = 0.004046916961669922
$ws.Range("K6").Value = 0.003969907760620117
$ws.Range("A7").Value = 1000
$ws.Range("B7").Value = 0.005184173583984375
$ws.Range("C7").Value = 0.004668235778808594
$ws.Range("D7").Value = 0.005610227584838867
$ws.Range("E7").Value = 0.004596471786499023
$ws.Range("F7").Value = 0.004528045654296875
$ws.Range("G7").Value = 0.004647016525268555
$ws.Range("H7").Value = 0.004456996917724609
$ws.Range("I7").Value = 0.00442957878112793
$ws.Range("J7").Value = 0.005900859832763672
$ws.Range("K7").Value = 0.005107402801513672

$ws = $wb.Worksheets.Item("Summary")
$ws.Range("A3").Value = 200
$ws.Range("A4").Value = 400
$ws.Range("A5").Value = 600
$ws.Range("A6").Value = 800
$ws.Range("A7").Value = 1000

